$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the dates between rows 2-3 and rows 4-5 in column D
$ws.Range("D2").Value = 44838
$ws.Range("D3").Value = 44838
$ws.Range("D4").Value = 44846
$ws.Range("D5").Value = 44846
